$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.026.85"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'  -0.41%  "
$ws.Range("E2").ClearFormats()
$ws.Range("D3").Value = "'2.420.45"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'  -0.14%  "
$ws.Range("E3").ClearFormats()
$ws.Range("E4").Value = "'  +0.04%  "
$ws.Range("E4").ClearFormats()
$ws.Range("D5").Value = "'561.90"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'  -0.24%  "
$ws.Range("E5").ClearFormats()
$ws.Range("D6").Value = "'143.54"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "'  -0.45%  "
$ws.Range("E6").ClearFormats()
$ws.Range("D8").Value = "'0.530"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "'  -0.44%  "
$ws.Range("E8").ClearFormats()
$ws.Range("D9").Value = "'2.419.07"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "'  -0.13%  "
$ws.Range("E9").ClearFormats()
$ws.Range("E10").Value = "'  -0.46%  "
$ws.Range("E10").ClearFormats()
$ws.Range("E11").Value = "'  +0.19%  "
$ws.Range("E11").ClearFormats()
$ws.Range("D12").Value = "'5.19"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "'  -3.27%  "
$ws.Range("E12").ClearFormats()
$ws.Range("E13").Value = "'  -1.35%  "
$ws.Range("E13").ClearFormats()
$ws.Range("D14").Value = "'26.17"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'  +0.98%  "
$ws.Range("E14").ClearFormats()
$ws.Range("E15").Value = "'  -2.34%  "
$ws.Range("E15").ClearFormats()
$ws.Range("D16").Value = "'2.857.55"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'  -0.15%  "
$ws.Range("E16").ClearFormats()
$ws.Range("D17").Value = "'61.933.94"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "'  -0.26%  "
$ws.Range("E17").ClearFormats()
$ws.Range("D18").Value = "'2.424.53"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "'  -0.03%  "
$ws.Range("E18").ClearFormats()
$ws.Range("E19").Value = "'  -0.70%  "
$ws.Range("E19").ClearFormats()
$ws.Range("D20").Value = "'322.73"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "'  -0.28%  "
$ws.Range("E20").ClearFormats()
$ws.Range("E21").Value = "'  -1.74%  "
$ws.Range("E21").ClearFormats()
$ws.Range("E22").Value = "'  +0.97%  "
$ws.Range("E22").ClearFormats()
$ws.Range("E23").Value = "'  +0.00%  "
$ws.Range("E23").ClearFormats()
$ws.Range("D24").Value = "'67.38"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "'  +2.68%  "
$ws.Range("E24").ClearFormats()
$ws.Range("E25").Value = "'  +1.74%  "
$ws.Range("E25").ClearFormats()
$ws.Range("D26").Value = "'8.69"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "'  -2.73%  "
$ws.Range("E26").ClearFormats()
$ws.Range("D27").Value = "'558.19"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "'  -4.25%  "
$ws.Range("E27").ClearFormats()
$ws.Range("D28").Value = "'2.539.32"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "'  -0.14%  "
$ws.Range("E28").ClearFormats()
$ws.Range("E29").Value = "'  +0.08%  "
$ws.Range("E29").ClearFormats()
$ws.Range("E30").Value = "'  -1.18%  "
$ws.Range("E30").ClearFormats()
$ws.Range("E31").Value = "'  -0.73%  "
$ws.Range("E31").ClearFormats()
$ws.Range("E32").Value = "'  -4.64%  "
$ws.Range("E32").ClearFormats()
$ws.Range("E33").Value = "'  -2.11%  "
$ws.Range("E33").ClearFormats()
$ws.Range("D34").Value = "'1.86"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "'  -1.14%  "
$ws.Range("E34").ClearFormats()
$ws.Range("E35").Value = "'  -2.99%  "
$ws.Range("E35").ClearFormats()
$ws.Range("E36").Value = "'  -0.06%  "
$ws.Range("E36").ClearFormats()
$ws.Range("E37").Value = "'  -1.05%  "
$ws.Range("E37").ClearFormats()
$ws.Range("E38").Value = "'  -1.28%  "
$ws.Range("E38").ClearFormats()
$ws.Range("E39").Value = "'  -4.38%  "
$ws.Range("E39").ClearFormats()
$ws.Range("D40").Value = "'152.47"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "'  -0.43%  "
$ws.Range("E40").ClearFormats()
$ws.Range("E41").Value = "'  -0.07%  "
$ws.Range("E41").ClearFormats()
$ws.Range("E42").Value = "'  -1.22%  "
$ws.Range("E42").ClearFormats()
$ws.Range("E43").Value = "'  +0.32%  "
$ws.Range("E43").ClearFormats()
$ws.Range("E44").Value = "'  -2.75%  "
$ws.Range("E44").ClearFormats()
$ws.Range("D45").Value = "'147.35"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "'  -1.94%  "
$ws.Range("E45").ClearFormats()
$ws.Range("D46").Value = "'3.64"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "'  -0.66%  "
$ws.Range("E46").ClearFormats()
$ws.Range("D47").Value = "'0.0528"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "'  -1.95%  "
$ws.Range("E47").ClearFormats()
$ws.Range("D48").Value = "'19.89"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "'  -2.15%  "
$ws.Range("E48").ClearFormats()
$ws.Range("E49").Value = "'  +0.00%  "
$ws.Range("E49").ClearFormats()
$ws.Range("D50").Value = "'0.0920"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "'  -0.49%  "
$ws.Range("E50").ClearFormats()
$ws.Range("E51").Value = "'  -0.17%  "
$ws.Range("E51").ClearFormats()
